# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the associated handoff timestamps, on all three
# sheets of the workbook (Overview, zh-cn, de-de). Also widens the
# "Status" columns so the longer "Ready for handoff" text fits.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# New column width (character units) for the COM object model. Excel
# rounds ColumnWidth to whole-pixel boundaries when it serializes the
# sheet, so this is the value that lands closest to the widened target
# width used by the report generator.
$newColumnWidth = 16.3

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-29 17:05:37"

$wsOverview.Columns("E").ColumnWidth = $newColumnWidth
$wsOverview.Columns("F").ColumnWidth = $newColumnWidth

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-29 17:05:33"

$wsZhCn.Columns("C").ColumnWidth = $newColumnWidth

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-29 17:05:37"

$wsDeDe.Columns("C").ColumnWidth = $newColumnWidth
